# Applies three textual edits to draft-gandhi-spring-stamp-srpm-01.pptx:
#  1. Slide 14 ("Example Provisioning Model"): change the dashed connector
#     between R1 and R5 boxes from "------------" to "============".
#  2. Slide 16 ("SRv6 STAMP / Probe Query Message"): collapse the
#     ".  <SID List>...." line's split runs into one merged run.
#  3. Slide 16: trim one trailing space from the "Reflector IPv6 Address"
#     line (the run that follows the bold "Reflector" run).

$p = $ppt.ActivePresentation

# --- Edit 1: slide 14, Rectangle 6 -----------------------------------
$s14 = $p.Slides.Item(14)
$sh14 = $s14.Shapes.Item(4)
$tr14 = $sh14.TextFrame.TextRange
$full14 = $tr14.Text
$needle14 = "                     |   R1  |------------|   R5  |"
$replacement14 = "                     |   R1  |============|   R5  |"
$idx14 = $full14.IndexOf($needle14)
if ($idx14 -ge 0) {
    $sub14 = $tr14.Characters($idx14 + 1, $needle14.Length)
    $sub14.Text = $replacement14
}

# --- Edit 2 & 3: slide 16, Rectangle 8 --------------------------------
$s16 = $p.Slides.Item(16)
$sh16 = $s16.Shapes.Item(6)
$tr16 = $sh16.TextFrame.TextRange

# Edit 2: merge ".  <" + "SID List" + ">   ..." + "." runs into one run.
$full16 = $tr16.Text
$needleSid = ".  <SID List>                                                   ."
$idxSid = $full16.IndexOf($needleSid)
if ($idxSid -ge 0) {
    $subSid = $tr16.Characters($idxSid + 1, $needleSid.Length)
    $subSid.Text = $needleSid
}

# Edit 3: drop one trailing space before the period in the run following
# the bold "Reflector" run.
$full16b = $tr16.Text
$needleReflector = "Reflector IPv6 Address               ."
$idxReflector = $full16b.IndexOf($needleReflector)
if ($idxReflector -ge 0) {
    $afterReflector = $idxReflector + "Reflector".Length
    $oldTail = " IPv6 Address               ."
    $newTail = " IPv6 Address              ."
    $subTail = $tr16.Characters($afterReflector + 1, $oldTail.Length)
    $subTail.Text = $newTail
}
